$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 5, shifting existing rows 5-8 down to 6-9.
# Excel will auto-adjust the SUM formulas in row 4 (G5:G9 -> G5:G10, etc.)
# and the relative AC formulas (SUM(S5,AB5) -> SUM(S6,AB6), etc.) as part
# of the insert.
$ws.Rows.Item(5).Insert()
# The inserted row picks up formatting copied from the row above (the
# "Totali" row); the data rows in this sheet carry no explicit cell
# styles, so strip the copied formatting to match.
$ws.Rows.Item(5).ClearFormats()

# Populate the newly inserted row 5 with the new invoice entry.
$ws.Cells.Item(5, 2).Value = "16.01.2024"
# Column C's text ("2024/01/16") looks like a date, so a plain .Value
# assignment would be auto-parsed into a date serial. Lead with an
# apostrophe (same as typing it directly into Excel) to force literal text.
$ws.Cells.Item(5, 3).Value = "'2024/01/16"
$ws.Cells.Item(5, 4).Value = "Magic Ice "
$ws.Cells.Item(5, 5).Value = 810379471
$ws.Cells.Item(5, 6).Value = 330065898
$ws.Cells.Item(5, 7).Value = 0
$ws.Cells.Item(5, 13).Value = 22.1
$ws.Cells.Item(5, 19).Value = 3.98
$ws.Cells.Item(5, 22).Value = 0
$ws.Cells.Item(5, 28).Value = 0
$ws.Range("AC5").Formula = "=SUM(S5,AB5)"

# Renumber the "Nr." column (A) for all data rows (now 1..5).
$ws.Cells.Item(5, 1).Value = 1
$ws.Cells.Item(6, 1).Value = 2
$ws.Cells.Item(7, 1).Value = 3
$ws.Cells.Item(8, 1).Value = 4
$ws.Cells.Item(9, 1).Value = 5

# Ensure the Totali row's SUM ranges cover the full data block (5:10),
# matching the target workbook exactly (the insert can otherwise anchor
# the range start one row too low).
$ws.Range("G4").Formula = "=SUM(G5:G10)"
$ws.Range("M4").Formula = "=SUM(M5:M10)"
$ws.Range("S4").Formula = "=SUM(S5:S10)"
$ws.Range("V4").Formula = "=SUM(V5:V10)"
$ws.Range("AB4").Formula = "=SUM(AB5:AB10)"
$ws.Range("AC4").Formula = "=SUM(AC5:AC10)"
